$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.923.78'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '2.462.21'
$ws.Range('E3').Value = '  -2.52%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''517.00'
$ws.Range('E5').Value = '  -3.64%  '
$ws.Range('D6').Value = '''130.76'
$ws.Range('E6').Value = '  -4.42%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').Value = '''0.557'
$ws.Range('E8').Value = '  -1.88%  '
$ws.Range('D9').Value = '2.464.10'
$ws.Range('E9').Value = '  -2.38%  '
$ws.Range('D10').Value = '''0.0987'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '''5.28'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('E13').Value = '  -2.56%  '
$ws.Range('D14').Value = '2.899.33'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '57.868.90'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '''22.25'
$ws.Range('E16').Value = '  -3.62%  '
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').Value = '2.460.64'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('D19').Value = '''10.68'
$ws.Range('E19').Value = '  -4.10%  '
$ws.Range('D20').Value = '''319.62'
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('E21').Value = '  -2.62%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '''5.69'
$ws.Range('E23').Value = '  -4.32%  '
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('E25').Value = '  -3.06%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -3.81%  '
$ws.Range('D28').Value = '''7.29'
$ws.Range('E28').Value = '  -3.08%  '
$ws.Range('D29').Value = '0.0₃0737'
$ws.Range('E29').Value = '  -4.50%  '
$ws.Range('D30').Value = '''165.52'
$ws.Range('E30').Value = '  -3.76%  '
$ws.Range('E31').Value = '  -4.69%  '
$ws.Range('D32').Value = '''6.23'
$ws.Range('E32').Value = '  -6.76%  '
$ws.Range('D33').Value = '''1.15'
$ws.Range('E33').Value = '  -2.00%  '
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').Value = '''17.99'
$ws.Range('E36').Value = '  -2.22%  '
$ws.Range('D37').Value = '''1.28'
$ws.Range('E37').Value = '  -9.29%  '
$ws.Range('D38').Value = '''3.94'
$ws.Range('E38').Value = '  -3.77%  '
$ws.Range('E39').Value = '  -5.17%  '
$ws.Range('D40').Value = '''0.784'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('D41').Value = '''3.44'
$ws.Range('D42').Value = '''269.82'
$ws.Range('E42').Value = '  -5.35%  '
$ws.Range('D43').Value = '''4.92'
$ws.Range('E43').Value = '  -5.38%  '
$ws.Range('E44').Value = '  -2.98%  '
$ws.Range('D45').Value = '''125.80'
$ws.Range('E45').Value = '  -4.28%  '
$ws.Range('D46').Value = '''0.0904'
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('D47').Value = '''0.0485'
$ws.Range('E47').Value = '  -4.40%  '
$ws.Range('D48').Value = '''0.0209'
$ws.Range('E48').Value = '  -5.11%  '
$ws.Range('D49').Value = '''16.77'
$ws.Range('E49').Value = '  -3.75%  '
$ws.Range('D50').Value = '1.715.52'
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('E51').Value = '  -2.10%  '
